$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 11

$ws.Cells.Item($newRow, 1).Value = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$ws.Cells.Item($newRow, 2).Value = "Beste klant,`nBedankt voor je interesse in de VentiQ-250. Helaas kunnen we op basis van je e-mailadres geen datasheet vinden. Zou je ons kunnen voorzien van meer informatie, zoals je volledige naam, bedrijfsnaam of eventuele andere gegevens waaronder de datasheet geregistreerd staat? Hiermee kunnen we je beter van dienst zijn en de datasheet naar je opsturen.`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent"
$ws.Cells.Item($newRow, 3).Value = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$ws.Cells.Item($newRow, 4).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 5).Value = "Productinformatie"
$ws.Cells.Item($newRow, 6).Value = "2025-07-31 21:50:21"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Nee"
$ws.Cells.Item($newRow, 9).Value = "Ja"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# The multi-line body text in column B would otherwise trigger an
# automatic custom row height; re-running AutoFit brings the row back
# to the sheet's standard height (matching the other data rows, which
# have no explicit row height despite also containing multi-line text).
$ws.Rows.Item($newRow).AutoFit()
